$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FEINmismatch")
$ws1.Range("B2").Value = "Thu Nov 07 15:55:03 EST 2024"
$ws1.Range("B3").Value = "Thu Nov 07 15:55:16 EST 2024"
$ws1.Range("B4").Value = "Thu Nov 07 15:55:28 EST 2024"
$ws1.Range("B5").Value = "Thu Nov 07 15:55:41 EST 2024"
$ws1.Range("B6").Value = "Thu Nov 07 15:55:54 EST 2024"
$ws1.Range("B7").Value = "Thu Nov 07 15:56:06 EST 2024"
$ws1.Range("B8").Value = "Thu Nov 07 15:56:18 EST 2024"
$ws1.Range("B9").Value = "Thu Nov 07 15:56:30 EST 2024"
$ws1.Range("B10").Value = "Thu Nov 07 15:56:42 EST 2024"
$ws1.Range("B11").Value = "Thu Nov 07 15:56:55 EST 2024"
$ws1.Range("B12").Value = "Thu Nov 07 15:57:07 EST 2024"
$ws1.Range("B13").Value = "Thu Nov 07 15:57:19 EST 2024"
$ws1.Range("B14").Value = "Thu Nov 07 15:57:33 EST 2024"
$ws1.Range("B15").Value = "Thu Nov 07 15:57:45 EST 2024"
$ws1.Range("B16").Value = "Thu Nov 07 15:57:57 EST 2024"
$ws1.Range("B17").Value = "Thu Nov 07 15:58:10 EST 2024"
$ws1.Range("B18").Value = "Thu Nov 07 15:58:22 EST 2024"
$ws1.Range("B19").Value = "Thu Nov 07 15:58:34 EST 2024"
$ws1.Range("B20").Value = "Thu Nov 07 15:58:47 EST 2024"
$ws1.Range("B21").Value = "Thu Nov 07 15:58:59 EST 2024"
$ws1.Range("B22").Value = "Thu Nov 07 15:59:13 EST 2024"
$ws1.Range("B23").Value = "Thu Nov 07 15:59:25 EST 2024"
$ws1.Range("B24").Value = "Thu Nov 07 15:59:38 EST 2024"
$ws1.Range("B25").Value = "Thu Nov 07 15:59:50 EST 2024"
$ws1.Range("B26").Value = "Thu Nov 07 16:00:03 EST 2024"
$ws1.Range("B27").Value = "Thu Nov 07 16:00:15 EST 2024"
$ws1.Range("B28").Value = "Thu Nov 07 16:00:27 EST 2024"
$ws1.Range("B29").Value = "Thu Nov 07 16:00:39 EST 2024"
$ws1.Range("B30").Value = "Thu Nov 07 16:00:52 EST 2024"

$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$ws2.Range("B2").Value = "Thu Nov 07 16:01:08 EST 2024"
$ws2.Range("B3").Value = "Thu Nov 07 16:01:20 EST 2024"
$ws2.Range("B4").Value = "Thu Nov 07 16:01:32 EST 2024"
$ws2.Range("B5").Value = "Thu Nov 07 16:01:44 EST 2024"
$ws2.Range("B6").Value = "Thu Nov 07 16:01:56 EST 2024"
$ws2.Range("B7").Value = "Thu Nov 07 16:02:08 EST 2024"
$ws2.Range("B8").Value = "Thu Nov 07 16:02:22 EST 2024"
$ws2.Range("B9").Value = "Thu Nov 07 16:02:35 EST 2024"
$ws2.Range("B10").Value = "Thu Nov 07 16:02:46 EST 2024"
$ws2.Range("B11").Value = "Thu Nov 07 16:02:59 EST 2024"
$ws2.Range("B12").Value = "Thu Nov 07 16:03:12 EST 2024"
$ws2.Range("B13").Value = "Thu Nov 07 16:03:24 EST 2024"
$ws2.Range("B14").Value = "Thu Nov 07 16:03:36 EST 2024"
$ws2.Range("B15").Value = "Thu Nov 07 16:03:47 EST 2024"
$ws2.Range("B16").Value = "Thu Nov 07 16:03:59 EST 2024"
$ws2.Range("B17").Value = "Thu Nov 07 16:04:12 EST 2024"
$ws2.Range("B18").Value = "Thu Nov 07 16:04:25 EST 2024"
$ws2.Range("B19").Value = "Thu Nov 07 16:04:55 EST 2024"
